$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 34-65 (values changed per diff)

$ws.Range("A34").Value = 32
$ws.Range("B34").Value = 12000001
$ws.Range("C34").Value = 14000034
$ws.Range("D34").Value = ""
$ws.Range("E34").Value = 'ITCH00001'

$ws.Range("A35").Value = 33
$ws.Range("B35").Value = 12000001
$ws.Range("C35").Value = 14000035
$ws.Range("D35").Value = ""
$ws.Range("E35").Value = 'ITCH00002'

$ws.Range("A36").Value = 34
$ws.Range("B36").Value = 12000001
$ws.Range("C36").Value = 14000038
$ws.Range("D36").Value = ""
$ws.Range("E36").Value = 'ITCH00004'

$ws.Range("A37").Value = 35
$ws.Range("B37").Value = 12000001
$ws.Range("C37").Value = 14000039
$ws.Range("D37").Value = ""
$ws.Range("E37").Value = 'ITCH00005'

$ws.Range("A38").Value = 36
$ws.Range("B38").Value = 12000001
$ws.Range("C38").Value = 14000040
$ws.Range("D38").Value = ""
$ws.Range("E38").Value = 'ITCH00006'

$ws.Range("A39").Value = 37
$ws.Range("B39").Value = 12000001
$ws.Range("C39").Value = 14000041
$ws.Range("D39").Value = ""
$ws.Range("E39").Value = 'ITCH00007'

$ws.Range("A40").Value = 38
$ws.Range("B40").Value = 12000001
$ws.Range("C40").Value = 14000046
$ws.Range("D40").Value = ""
$ws.Range("E40").Value = 'ITCH00011'

$ws.Range("A41").Value = 39
$ws.Range("B41").Value = 12000001
$ws.Range("C41").Value = 14000049
$ws.Range("D41").Value = ""
$ws.Range("E41").Value = 'ITCH00013'

$ws.Range("A42").Value = 40
$ws.Range("B42").Value = 12000001
$ws.Range("C42").Value = 14000050
$ws.Range("D42").Value = ""
$ws.Range("E42").Value = 'ITCH00014'

$ws.Range("A43").Value = 41
$ws.Range("B43").Value = 12000001
$ws.Range("C43").Value = 14000051
$ws.Range("D43").Value = ""
$ws.Range("E43").Value = 'ITCH00015'

$ws.Range("A44").Value = 42
$ws.Range("B44").Value = 12000001
$ws.Range("C44").Value = 14000052
$ws.Range("D44").Value = ""
$ws.Range("E44").Value = 'ITCH00016'

$ws.Range("A45").Value = 43
$ws.Range("B45").Value = 12000001
$ws.Range("C45").Value = 14000053
$ws.Range("D45").Value = ""
$ws.Range("E45").Value = 'ITCH00017'

$ws.Range("A46").Value = 44
$ws.Range("B46").Value = 12000001
$ws.Range("C46").Value = 14000054
$ws.Range("D46").Value = ""
$ws.Range("E46").Value = 'ITCH00018'

$ws.Range("A47").Value = 45
$ws.Range("B47").Value = 12000001
$ws.Range("C47").Value = 14000055
$ws.Range("D47").Value = ""
$ws.Range("E47").Value = 'ITCH00019'

$ws.Range("A48").Value = 46
$ws.Range("B48").Value = 12000001
$ws.Range("C48").Value = 14000056
$ws.Range("D48").Value = ""
$ws.Range("E48").Value = 'ITCH00020'

$ws.Range("A49").Value = 47
$ws.Range("B49").Value = 12000001
$ws.Range("C49").Value = 14000057
$ws.Range("D49").Value = ""
$ws.Range("E49").Value = 'ITCH00021'

$ws.Range("A50").Value = 48
$ws.Range("B50").Value = 12000001
$ws.Range("C50").Value = 14000058
$ws.Range("D50").Value = ""
$ws.Range("E50").Value = 'ITCH00022'

$ws.Range("A51").Value = 49
$ws.Range("B51").Value = 12000001
$ws.Range("C51").Value = 14000059
$ws.Range("D51").Value = ""
$ws.Range("E51").Value = 'ITCH00023'

$ws.Range("A52").Value = 50
$ws.Range("B52").Value = 12000001
$ws.Range("C52").Value = 14000060
$ws.Range("D52").Value = ""
$ws.Range("E52").Value = 'ITCH00024'

$ws.Range("A53").Value = 51
$ws.Range("B53").Value = 12000001
$ws.Range("C53").Value = 14000063
$ws.Range("D53").Value = ""
$ws.Range("E53").Value = 'ITCH00026'

$ws.Range("A54").Value = 52
$ws.Range("B54").Value = 12000001
$ws.Range("C54").Value = 14000066
$ws.Range("D54").Value = ""
$ws.Range("E54").Value = 'ITCH00028'

$ws.Range("A55").Value = 53
$ws.Range("B55").Value = 12000001
$ws.Range("C55").Value = 14000068
$ws.Range("D55").Value = ""
$ws.Range("E55").Value = 'ITCH00029'

$ws.Range("A56").Value = 54
$ws.Range("B56").Value = 12000001
$ws.Range("C56").Value = 14000071
$ws.Range("D56").Value = ""
$ws.Range("E56").Value = 'ITCH00032'

$ws.Range("A57").Value = 55
$ws.Range("B57").Value = 12000001
$ws.Range("C57").Value = 14000072
$ws.Range("D57").Value = ""
$ws.Range("E57").Value = 'ITCH00033'

$ws.Range("A58").Value = 56
$ws.Range("B58").Value = 12000001
$ws.Range("C58").Value = 14000073
$ws.Range("D58").Value = ""
$ws.Range("E58").Value = 'ITCH00034'

$ws.Range("A59").Value = 57
$ws.Range("B59").Value = 12000001
$ws.Range("C59").Value = 14000074
$ws.Range("D59").Value = ""
$ws.Range("E59").Value = 'ITCH00035'

$ws.Range("A60").Value = 58
$ws.Range("B60").Value = 12000001
$ws.Range("C60").Value = 14000075
$ws.Range("D60").Value = ""
$ws.Range("E60").Value = 'ITCH00036'

$ws.Range("A61").Value = 59
$ws.Range("B61").Value = 12000001
$ws.Range("C61").Value = 14000076
$ws.Range("D61").Value = ""
$ws.Range("E61").Value = 'ITCH00040'

$ws.Range("A62").Value = 60
$ws.Range("B62").Value = 12000002
$ws.Range("C62").Value = 14000069
$ws.Range("D62").Value = 'Definizione italiana, oggetto svizzero'
$ws.Range("E62").Value = 'ITCH00030'

$ws.Range("A63").Value = 61
$ws.Range("B63").Value = 12000036
$ws.Range("C63").Value = 14000036
$ws.Range("D63").Value = 'Definizioni riorganizzate da Tommaso Sansone, Politecnico di Milano.'
$ws.Range("E63").Value = 'ITCH00003'

$ws.Range("A64").Value = 62
$ws.Range("B64").Value = 12000042
$ws.Range("C64").Value = 14000042
$ws.Range("D64").Value = 'Per il lato italiano, definizione tratta dall''art 10 del codice di protezione civile del 2-1-2018, rielaborata da Tommaso Sansone.'
$ws.Range("E64").Value = 'ITCH00009'

$ws.Range("A65").Value = 63
$ws.Range("B65").Value = 12000044
$ws.Range("C65").Value = 14000044
$ws.Range("D65").Value = 'strutture operative in italia e organizzazioni partner in svizzera hanno la stessa connotazione all''interno di un sistema di protezione della popolazione.'
$ws.Range("E65").Value = 'ITCH00010'

# Rows 66-69 are brand new; copy the header-row-A style (bold/border/centered)
# from an existing styled A-column cell so formatting matches the rest of column A.
$ws.Range("A65").Copy()
$ws.Range("A66:A69").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A66").Value = 64
$ws.Range("B66").Value = 12000047
$ws.Range("C66").Value = 14000047
$ws.Range("D66").Value = 'Per il lato italiano, definizione rielaborata da Tommaso Sansone a partire dall''art 7 del decreto di protezione civile 2-1-2018'
$ws.Range("E66").Value = 'ITCH00012'

$ws.Range("A67").Value = 65
$ws.Range("B67").Value = 12000061
$ws.Range("C67").Value = 14000061
$ws.Range("D67").Value = 'Il significato italiano di protezione civile coincide a livello di strutture con io significato svizzero di protezione della popolazione.' + "`n" + 'Definizioni riorganizzate da Tommaso Sansone, Politecnico di Milano.'
$ws.Range("E67").Value = 'ITCH00025'

$ws.Range("A68").Value = 66
$ws.Range("B68").Value = 12000064
$ws.Range("C68").Value = 14000064
$ws.Range("D68").Value = 'definizione lato italiano, oggetto lato svizzero'
$ws.Range("E68").Value = 'ITCH00027'

$ws.Range("A69").Value = 67
$ws.Range("B69").Value = 12000070
$ws.Range("C69").Value = 14000070
$ws.Range("D69").Value = 'definizione italiana, oggetto svizzero'
$ws.Range("E69").Value = 'ITCH00031'
